$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.517.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.874.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +1.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'313.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.59%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.014"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.4788"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3773"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.69%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07373"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.50%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9385"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.44%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +5.60%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07857"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.24%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.903.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.44%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.450"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.72%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.594"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.07%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'90.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.73%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.017"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.91%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008931"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.34%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.94%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.75%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'27.566.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.45%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.140"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.15%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.04%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.961"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.17%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'153.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.07%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'18.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.46%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.019"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.94%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'115.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.46%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.003"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.28%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.08937"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.02%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.328"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.81%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +3.72%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.614"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.00%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7495"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.29%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +6.73%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D38").Value = "'0.05305"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.87%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.007"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.59%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.5344"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.67%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'7.103"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.70%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1527"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.96%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.427"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.67%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'10.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.03%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.4837"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.00%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.015"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "'  +3.48%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'103.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.69%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'67.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.90%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06098"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.8993"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.80%  "
$ws.Range("E51").Style = "Normal"
